$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New table data (player, position, team) to replace the existing rows 2-17.
# Net effect vs. original: the "Ausar Thompson | SF,PF | Detroit Pistons" row is
# removed and the remaining rows are reordered.
$data = @(
    @("Tyrese Maxey", "PG,SG", "Philadelphia 76ers"),
    @("Austin Reaves", "PG,SG", "Los Angeles Lakers"),
    @("Stephen Curry", "PG,SG", "Golden State Warriors"),
    @("Tyrese Haliburton", "PG,SG", "Indiana Pacers"),
    @("Franz Wagner", "SF,PF", "Orlando Magic"),
    @("OG Anunoby", "SF,PF", "New York Knicks"),
    @("Kevin Durant", "SF,PF", "Phoenix Suns"),
    @("Keegan Murray", "SF,PF", "Sacramento Kings"),
    @("Mark Williams", "C", "Charlotte Hornets"),
    @("Karl-Anthony Towns", "PF,C", "New York Knicks"),
    @("Jarrett Allen", "C", "Cleveland Cavaliers"),
    @("Trey Murphy III", "SF,PF", "New Orleans Pelicans"),
    @("Jalen Duren", "C", "Detroit Pistons"),
    @("Daniel Gafford", "PF,C", "Dallas Mavericks"),
    @("Darius Garland", "PG", "Cleveland Cavaliers")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

# The old table had one extra row (row 17); remove it entirely so the sheet
# dimension shrinks back to A1:C16.
$lastRow = 17
$ws.Range($ws.Cells.Item($lastRow, 1), $ws.Cells.Item($lastRow, 3)).Delete()
